$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data per latest scrape run.
# For D-column cells whose new value looks like a plain number (one
# decimal point, e.g. "86.94"), set NumberFormat to Text ("@") first so
# Excel stores the literal digit-dot-digit string rather than silently
# parsing it into a numeric cell -- matching the source data which is
# plain text (note some prices use "." as a thousands separator too,
# e.g. "39.770.86", which already can't parse as a number).

    $ws.Range("D2").Value = "39.770.86"
    $ws.Range("E2").Value = "  +0.11%  "
    $ws.Range("D3").Value = "2.214.58"
    $ws.Range("E3").Value = "  +0.49%  "
    $ws.Range("E4").Value = "  -0.05%  "
    $ws.Range("D5").NumberFormat = "@"
    $ws.Range("D5").Value = "291.19"
    $ws.Range("E5").Value = "  -0.21%  "
    $ws.Range("D6").NumberFormat = "@"
    $ws.Range("D6").Value = "86.94"
    $ws.Range("E6").Value = "  +0.62%  "
    $ws.Range("D7").NumberFormat = "@"
    $ws.Range("D7").Value = "0.512"
    $ws.Range("E7").Value = "  -0.39%  "
    $ws.Range("E8").Value = "  +0.01%  "
    $ws.Range("D9").NumberFormat = "@"
    $ws.Range("D9").Value = "0.466"
    $ws.Range("E9").Value = "  -0.88%  "
    $ws.Range("D10").NumberFormat = "@"
    $ws.Range("D10").Value = "30.35"
    $ws.Range("E10").Value = "  +0.38%  "
    $ws.Range("D11").NumberFormat = "@"
    $ws.Range("D11").Value = "0.0780"
    $ws.Range("E11").Value = "  -0.30%  "
    $ws.Range("D12").NumberFormat = "@"
    $ws.Range("D12").Value = "49.88"
    $ws.Range("E12").Value = "  +5.36%  "
    $ws.Range("E13").Value = "  +2.48%  "
    $ws.Range("D14").NumberFormat = "@"
    $ws.Range("D14").Value = "6.44"
    $ws.Range("E14").Value = "  +2.19%  "
    $ws.Range("D15").Value = "2.563.57"
    $ws.Range("E15").Value = "  +0.53%  "
    $ws.Range("D16").NumberFormat = "@"
    $ws.Range("D16").Value = "13.77"
    $ws.Range("E16").Value = "  -1.30%  "
    $ws.Range("D17").Value = "2.256.30"
    $ws.Range("E17").Value = "  +1.59%  "
    $ws.Range("D18").NumberFormat = "@"
    $ws.Range("D18").Value = "0.729"
    $ws.Range("E18").Value = "  +0.43%  "
    $ws.Range("D19").Value = "39.755.19"
    $ws.Range("E19").Value = "  +0.23%  "
    $ws.Range("D20").Value = "0.0₃0885"
    $ws.Range("E20").Value = "  +0.87%  "
    $ws.Range("D21").NumberFormat = "@"
    $ws.Range("D21").Value = "11.07"
    $ws.Range("E21").Value = "  -1.83%  "
    $ws.Range("D22").NumberFormat = "@"
    $ws.Range("D22").Value = "5.73"
    $ws.Range("E22").Value = "  -0.73%  "
    $ws.Range("D23").NumberFormat = "@"
    $ws.Range("D23").Value = "65.57"
    $ws.Range("E23").Value = "  -0.01%  "
    $ws.Range("D24").NumberFormat = "@"
    $ws.Range("D24").Value = "237.03"
    $ws.Range("E24").Value = "  +0.81%  "
    $ws.Range("E25").Value = "  +0.10%  "
    $ws.Range("D26").NumberFormat = "@"
    $ws.Range("D26").Value = "2.44"
    $ws.Range("E26").Value = "  -0.51%  "
    $ws.Range("E27").Value = "  +0.10%  "
    $ws.Range("D28").NumberFormat = "@"
    $ws.Range("D28").Value = "22.97"
    $ws.Range("E28").Value = "  +1.42%  "
    $ws.Range("B29").Value = "Toncoin"
    $ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    $ws.Range("D29").NumberFormat = "@"
    $ws.Range("D29").Value = "2.15"
    $ws.Range("E29").Value = "  -2.17%  "
    $ws.Range("B30").Value = "Cosmos"
    $ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    $ws.Range("D30").NumberFormat = "@"
    $ws.Range("D30").Value = "9.21"
    $ws.Range("E30").Value = "  -0.22%  "
    $ws.Range("D31").NumberFormat = "@"
    $ws.Range("D31").Value = "156.46"
    $ws.Range("E31").Value = "  +3.07%  "
    $ws.Range("D32").NumberFormat = "@"
    $ws.Range("D32").Value = "31.84"
    $ws.Range("E32").Value = "  -2.10%  "
    $ws.Range("D33").NumberFormat = "@"
    $ws.Range("D33").Value = "1.00"
    $ws.Range("E33").Value = "  +0.07%  "
    $ws.Range("D34").NumberFormat = "@"
    $ws.Range("D34").Value = "4.94"
    $ws.Range("E34").Value = "  +0.84%  "
    $ws.Range("D35").NumberFormat = "@"
    $ws.Range("D35").Value = "2.96"
    $ws.Range("E35").Value = "  +6.73%  "
    $ws.Range("D36").NumberFormat = "@"
    $ws.Range("D36").Value = "0.0711"
    $ws.Range("E36").Value = "  -0.78%  "
    $ws.Range("E37").Value = "  -1.77%  "
    $ws.Range("E38").Value = "  -0.34%  "
    $ws.Range("D39").NumberFormat = "@"
    $ws.Range("D39").Value = "0.0988"
    $ws.Range("E39").Value = "  +0.54%  "
    $ws.Range("D40").NumberFormat = "@"
    $ws.Range("D40").Value = "1.73"
    $ws.Range("E40").Value = "  +2.49%  "
    $ws.Range("D41").NumberFormat = "@"
    $ws.Range("D41").Value = "15.28"
    $ws.Range("E41").Value = "  -3.75%  "
    $ws.Range("D42").Value = "2.106.87"
    $ws.Range("E42").Value = "  +2.36%  "
    $ws.Range("D43").NumberFormat = "@"
    $ws.Range("D43").Value = "3.72"
    $ws.Range("E43").Value = "  -1.24%  "
    $ws.Range("D44").NumberFormat = "@"
    $ws.Range("D44").Value = "18.19"
    $ws.Range("E44").Value = "  +2.21%  "
    $ws.Range("E45").Value = "  +0.94%  "
    $ws.Range("D46").NumberFormat = "@"
    $ws.Range("D46").Value = "9.93"
    $ws.Range("E46").Value = "  +0.93%  "
    $ws.Range("D47").NumberFormat = "@"
    $ws.Range("D47").Value = "2.01"
    $ws.Range("E47").Value = "  -5.37%  "
    $ws.Range("D48").NumberFormat = "@"
    $ws.Range("D48").Value = "2.70"
    $ws.Range("E48").Value = "  +3.85%  "
    $ws.Range("D49").Value = "2.436.41"
    $ws.Range("E49").Value = "  +0.29%  "
    $ws.Range("D50").NumberFormat = "@"
    $ws.Range("D50").Value = "1.46"
    $ws.Range("E50").Value = "  +2.84%  "
    $ws.Range("E51").Value = "  +2.45%  "
